$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 84
$ws.Range("I33").Value = 107.1
$ws.Range("K33").Value = 107.1
$ws.Range("M33").Value = 121.9

$ws.Range("H40").Value = 5703.1665
$ws.Range("I40").Value = 4271.222
$ws.Range("K40").Value = 4271.222
$ws.Range("M40").Value = -4096.222

$ws.Range("H51").Value = 4500
$ws.Range("J51").Value = 5000
$ws.Range("L51").Value = 5000
$ws.Range("N51").Value = -5968

$ws.Range("H64").Value = 9499.25
$ws.Range("I64").Value = 8998.5
$ws.Range("K64").Value = 8998.5
$ws.Range("M64").Value = -8750.5

$ws.Range("H67").Value = 9499.25
$ws.Range("I67").Value = 8998.5
$ws.Range("K67").Value = 8998.5
$ws.Range("M67").Value = -8140.5

$ws.Range("H69").Value = 6901.6763
$ws.Range("J69").Value = 6899.0312
$ws.Range("L69").Value = 20697.0936
$ws.Range("N69").Value = -22445.0936

$ws.Range("H72").Value = 6901.6763
$ws.Range("J72").Value = 6899.0312
$ws.Range("L72").Value = 62091.2808
$ws.Range("N72").Value = -70827.28080000001

$ws.Range("H86").Value = 2901
$ws.Range("I86").Value = 2901
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 2901
$ws.Range("L86").Value = 0
$ws.Range("N86").Value = -1778

$ws.Range("H89").Value = 2901
$ws.Range("I89").Value = 2901
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 14505
$ws.Range("L89").Value = 0
$ws.Range("N89").Value = -8889

$ws.Range("H129").Value = 1584.6316
$ws.Range("I129").Value = 1015.4167
$ws.Range("J129").Value = 2560.4285
$ws.Range("K129").Value = 3046.2501
$ws.Range("L129").Value = 7681.2855
$ws.Range("M129").Value = 1953.7499
$ws.Range("N129").Value = -17681.2855

$ws.Range("H132").Value = 11719.087
$ws.Range("I132").Value = 12597.588
$ws.Range("K132").Value = 37792.764
$ws.Range("M132").Value = -35262.764

$ws.Range("H138").Value = 1680
$ws.Range("I138").Value = 1680
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 5040
$ws.Range("L138").Value = 0
$ws.Range("N138").Value = 100

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2746.75
$ws.Range("I45").Value = 1995.7778
$ws.Range("J45").Value = 4999.6665
$ws.Range("K45").Value = 1995.7778
$ws.Range("L45").Value = 4999.6665
$ws.Range("M45").Value = -1618.7778
$ws.Range("N45").Value = -5753.6665

$ws.Range("H74").Value = 4569.2
$ws.Range("I74").Value = 3229.25
$ws.Range("K74").Value = 3229.25
$ws.Range("M74").Value = -2355.25

$ws.Range("H77").Value = 4569.2
$ws.Range("I77").Value = 3229.25
$ws.Range("K77").Value = 16146.25
$ws.Range("M77").Value = -11778.25

$ws.Range("H88").Value = 1253.7273
$ws.Range("I88").Value = 954.8
$ws.Range("J88").Value = 1502.8334
$ws.Range("K88").Value = 954.8
$ws.Range("L88").Value = 1502.8334
$ws.Range("M88").Value = -548.8
$ws.Range("N88").Value = -2314.8334

$ws.Range("H91").Value = 1253.7273
$ws.Range("I91").Value = 954.8
$ws.Range("J91").Value = 1502.8334
$ws.Range("K91").Value = 954.8
$ws.Range("L91").Value = 1502.8334
$ws.Range("M91").Value = 449.2
$ws.Range("N91").Value = -4310.8334

$ws.Range("H110").Value = 3175.7693
$ws.Range("I110").Value = 2773.75
$ws.Range("K110").Value = 2773.75
$ws.Range("M110").Value = -728.75

$ws.Range("H122").Value = 1510
$ws.Range("I122").Value = 1510
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 4530
$ws.Range("L122").Value = 0
$ws.Range("N122").Value = -2080

$ws.Range("H132").Value = 2655.3333
$ws.Range("I132").Value = 1987.375
$ws.Range("K132").Value = 5962.125
$ws.Range("M132").Value = -3432.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4902.357
$ws.Range("I86").Value = 1184.1428
$ws.Range("J86").Value = 8620.571
$ws.Range("K86").Value = 1184.1428
$ws.Range("L86").Value = 8620.571
$ws.Range("M86").Value = -61.14280000000008
$ws.Range("N86").Value = -10866.571

$ws.Range("H89").Value = 4902.357
$ws.Range("I89").Value = 1184.1428
$ws.Range("J89").Value = 8620.571
$ws.Range("K89").Value = 5920.714
$ws.Range("L89").Value = 43102.855
$ws.Range("M89").Value = -304.7139999999999
$ws.Range("N89").Value = -54334.855

$ws.Range("H94").Value = 1664.3334
$ws.Range("I94").Value = 1664.3334
$ws.Range("K94").Value = 1664.3334
$ws.Range("M94").Value = -1213.3334

$ws.Range("H99").Value = 3720.5
$ws.Range("I99").Value = 3313.9
$ws.Range("K99").Value = 3313.9
$ws.Range("M99").Value = -1815.9

$ws.Range("H107").Value = 5173.773
$ws.Range("I107").Value = 1082.6
$ws.Range("K107").Value = 1082.6
$ws.Range("M107").Value = 837.4000000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 3979.1667
$ws.Range("I22").Value = 3458.3333
$ws.Range("K22").Value = 3458.3333
$ws.Range("M22").Value = -3108.3333

$ws.Range("H60").Value = 21551.5
$ws.Range("I60").Value = 0
$ws.Range("K60").Value = 0
$ws.Range("M60").Value = ""

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 964
$ws.Range("J32").Value = 983
$ws.Range("L32").Value = 2949
$ws.Range("N32").Value = -3515

$ws.Range("H109").Value = 93496
$ws.Range("I109").Value = 126869.625
$ws.Range("J109").Value = 4499.6665
$ws.Range("K109").Value = 380608.875
$ws.Range("L109").Value = 13498.9995
$ws.Range("M109").Value = -379568.875
$ws.Range("N109").Value = -15578.9995

$ws.Range("H114").Value = 1770.75
$ws.Range("I114").Value = 1382.75
$ws.Range("J114").Value = 1964.75
$ws.Range("K114").Value = 4148.25
$ws.Range("L114").Value = 5894.25
$ws.Range("M114").Value = -894.25
$ws.Range("N114").Value = -12402.25

$ws.Range("H131").Value = 2162.0527
$ws.Range("I131").Value = 1858
$ws.Range("K131").Value = 5574
$ws.Range("M131").Value = -534

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6496.8335
$ws.Range("I70").Value = 5796.4
$ws.Range("K70").Value = 5796.4
$ws.Range("M70").Value = -5526.4

$ws.Range("H73").Value = 6496.8335
$ws.Range("I73").Value = 5796.4
$ws.Range("K73").Value = 5796.4
$ws.Range("M73").Value = -4860.4

$ws.Range("H80").Value = 2240.8333
$ws.Range("I80").Value = 1650
$ws.Range("J80").Value = 2831.6667
$ws.Range("K80").Value = 1650
$ws.Range("L80").Value = 2831.6667
$ws.Range("M80").Value = -652
$ws.Range("N80").Value = -4827.6667

$ws.Range("H83").Value = 2240.8333
$ws.Range("I83").Value = 1650
$ws.Range("J83").Value = 2831.6667
$ws.Range("K83").Value = 8250
$ws.Range("L83").Value = 14158.3335
$ws.Range("M83").Value = -3258
$ws.Range("N83").Value = -24142.3335

$ws.Range("H105").Value = 27650
$ws.Range("J105").Value = 27650
$ws.Range("L105").Value = 27650
$ws.Range("N105").Value = -34638

$ws.Range("H113").Value = 8044.4546
$ws.Range("I113").Value = 4623.5
$ws.Range("K113").Value = 4623.5
$ws.Range("M113").Value = -2453.5

$ws.Range("H122").Value = 336820.47
$ws.Range("I122").Value = 558660.25
$ws.Range("J122").Value = 4060.8333
$ws.Range("K122").Value = 1675980.75
$ws.Range("L122").Value = 12182.4999
$ws.Range("M122").Value = -1673530.75
$ws.Range("N122").Value = -17082.4999

$ws.Range("H132").Value = 83205.38
$ws.Range("I132").Value = 96515.55
$ws.Range("K132").Value = 289546.65
$ws.Range("M132").Value = -287016.65

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 5789
$ws.Range("J22").Value = 5789
$ws.Range("L22").Value = 5789
$ws.Range("N22").Value = -6379

$ws.Range("H27").Value = 5789
$ws.Range("J27").Value = 5789
$ws.Range("L27").Value = 5789
$ws.Range("N27").Value = -6003

$ws.Range("H40").Value = 7576.8887
$ws.Range("I40").Value = 6490.6665
$ws.Range("K40").Value = 6490.6665
$ws.Range("M40").Value = -6354.6665

$ws.Range("H61").Value = 6598
$ws.Range("I61").Value = 5596.5713
$ws.Range("K61").Value = 5596.5713
$ws.Range("M61").Value = -5394.5713

$ws.Range("H93").Value = 1628.3636

$ws.Range("H113").Value = 6598
$ws.Range("I113").Value = 5596.5713
$ws.Range("K113").Value = 5596.5713
$ws.Range("M113").Value = -3426.5713

$ws.Range("H136").Value = 3982.2
$ws.Range("I136").Value = 1305.5
$ws.Range("K136").Value = 3916.5
$ws.Range("M136").Value = -1366.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H105").Value = 70615
$ws.Range("J105").Value = 70615
$ws.Range("L105").Value = 70615
$ws.Range("N105").Value = -77603
